$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Range("B4").Value = 20
$ws.Range("C4").Value = 30
$ws.Range("D4").Value = 50
$ws.Range("F4").Value = "Pass"
$ws.Range("G4").Value = "Hibás"

# Row 5
$ws.Range("B5").Value = 30
$ws.Range("C5").Value = 70
$ws.Range("D5").Value = 100
$ws.Range("F5").Value = "Pass"

# Row 6
$ws.Range("B6").Value = 90
$ws.Range("C6").Value = 60
$ws.Range("D6").Value = 150
$ws.Range("F6").Value = "Pass"
$ws.Range("G6").Value = "Hibás"

# Row 7 stays empty, but gets touched (group/ungroup) leaving a bare row marker
$ws.Rows.Item(7).Group() | Out-Null
$ws.Rows.Item(7).Ungroup() | Out-Null

# Far below rows, only the "Üzenet" column populated
$ws.Range("G34").Value = "Hibás"
$ws.Range("G36").Value = "Hibás"

# Final selection highlights the freshly entered input cells
$ws.Range("B3:D6").Select() | Out-Null
